# SCD0025-003 update:
#  - Rename sheet SCD0333 -> SCD0025
#  - Update TC_ID (B2) from "DGS-348" to "SCD0025-003"
#  - Column B widens (bestFit) to fit the new, longer TC_ID text
#  - Active selection moves to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "SCD0025"

# Update the TC_ID cell value
$ws.Range("B2").Value = "SCD0025-003"

# Resize column B to (best-)fit the new text
$ws.Columns("B").ColumnWidth = 11.6

# Move the active selection to B3, matching the saved view state
$ws.Range("B3").Select() | Out-Null

Write-Host "Updated sheet name, TC_ID, column width and selection."
